$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cost values in column B (rows 3-7)
$ws.Range("B3").Value = 29550
$ws.Range("B4").Value = 27850
$ws.Range("B5").Value = 24850
$ws.Range("B6").Value = 22750
$ws.Range("B7").Value = 15150

# Remove rows 8 and 9 (time=35 and time=40 data points no longer used)
$ws.Range("A8:C9").ClearContents()

# Update the chart series so they only reference rows 2-7
$chart = $ws.ChartObjects().Item(1).Chart
$series1 = $chart.SeriesCollection().Item(1)
$series1.Formula = "=SERIES(Kosten!`$B`$1,Kosten!`$A`$2:`$A`$7,Kosten!`$B`$2:`$B`$7,1)"

$series2 = $chart.SeriesCollection().Item(2)
$series2.Formula = "=SERIES(Kosten!`$C`$1,Kosten!`$A`$2:`$A`$7,Kosten!`$C`$2:`$C`$7,2)"
